# fix contador materias completadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing "Contabilidad" row data: date and grade
$ws.Range("C6").Value = [DateTime]"2016-09-24"
$ws.Range("D6").Value = 9

# Fix the "Completadas" (completed count) formula to count all three
# grade columns (D, G, J), not just D
$ws.Range("C19").Formula = "=+COUNTA(D2:D14,G2:G14,J2:J14)"

# Move the selection to D2
$ws.Range("D2").Select()

$wb.Save()
